$d = $word.ActiveDocument

# The site rebuild dropped the trailing "Ver no Jupiter / Salvar em pdf / Salvar
# em docx" line, the "(c) 2020 ..." footer line, and the blank paragraph that
# separated them from the bibliography entry above. Locate those three
# paragraphs by their text (robust to their exact index) and delete the
# contiguous range they occupy, leaving the bibliography entry and the
# trailing blank paragraph (before the page break) untouched.

$jupiterText = "Ver no Jupiter Salvar em pdf Salvar em docx"
$copyrightText = [char]0xA9 + " 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution"

$paraCount = $d.Paragraphs.Count

$jupiterIdx = -1
$copyrightIdx = -1

for ($i = 1; $i -le $paraCount; $i++) {
    $text = $d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13)
    if ($text -eq $jupiterText) {
        $jupiterIdx = $i
    }
    if ($text -eq $copyrightText) {
        $copyrightIdx = $i
    }
}

if ($jupiterIdx -gt 0 -and $copyrightIdx -eq ($jupiterIdx + 1)) {
    # Include the blank paragraph immediately preceding the "Ver no Jupiter"
    # line so it is removed along with the two text paragraphs.
    $startIdx = $jupiterIdx - 1

    $startPara = $d.Paragraphs.Item($startIdx)
    $endPara = $d.Paragraphs.Item($copyrightIdx)

    $deleteRange = $d.Range($startPara.Range.Start, $endPara.Range.End)
    $deleteRange.Delete()
}
